$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5856198072433472
$ws.Range("B1").Value = 1.101048588752747
$ws.Range("C1").Value = 5.315079212188721
$ws.Range("D1").Value = 1.85684597492218
$ws.Range("E1").Value = 1.040079236030579
